$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $r = $ws.Range($cellRef)
    $r.NumberFormat = "@"
    $r.Value = $value
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "29.469.45"
$ws.Range("E2").Value = "  -1.04%  "
$ws.Range("D3").Value = "1.849.78"
$ws.Range("E3").Value = "  -0.56%  "
Set-TextValue "D4" "0.9990"
$ws.Range("E4").Value = "  -0.08%  "
Set-TextValue "D5" "241.81"
$ws.Range("E5").Value = "  -1.34%  "
Set-TextValue "D6" "0.6256"
$ws.Range("E6").Value = "  -2.58%  "
Set-TextValue "D7" "1.000"
Set-TextValue "D8" "48.22"
$ws.Range("E8").Value = "  +1.41%  "
Set-TextValue "D9" "0.07527"
$ws.Range("E9").Value = "  -0.19%  "
Set-TextValue "D10" "0.2975"
$ws.Range("E10").Value = "  -0.03%  "
$ws.Range("E11").Value = "  -1.43%  "
$ws.Range("D12").Value = "1.977.10"
$ws.Range("E12").Value = "  +6.29%  "
Set-TextValue "D13" "0.07699"
$ws.Range("E13").Value = "  +0.31%  "
Set-TextValue "D14" "5.001"
$ws.Range("E14").Value = "  -0.93%  "
Set-TextValue "D15" "0.6848"
$ws.Range("E15").Value = "  -1.23%  "
Set-TextValue "D16" "83.77"
$ws.Range("E16").Value = "  -0.31%  "
Set-TextValue "D17" "0.000009736"
$ws.Range("E17").Value = "  -1.08%  "
$ws.Range("D18").Value = "2.214.65"
$ws.Range("E18").Value = "  +4.78%  "
Set-TextValue "D19" "6.221"
$ws.Range("E19").Value = "  +1.46%  "
$ws.Range("D20").Value = "29.618.10"
$ws.Range("E20").Value = "  -0.57%  "
Set-TextValue "D21" "233.91"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("E22").Value = "  -1.38%  "
Set-TextValue "D23" "0.9999"
$ws.Range("E23").Value = "  -0.01%  "
Set-TextValue "D24" "7.590"
$ws.Range("E24").Value = "  +0.77%  "
Set-TextValue "D25" "1.0000"
$ws.Range("E25").Value = "  -0.10%  "
Set-TextValue "D26" "155.25"
$ws.Range("E26").Value = "  -2.31%  "
Set-TextValue "D27" "0.1389"
$ws.Range("E27").Value = "  -2.50%  "
Set-TextValue "D28" "8.436"
$ws.Range("E28").Value = "  -1.50%  "
Set-TextValue "D29" "17.71"
$ws.Range("E29").Value = "  -1.29%  "
Set-TextValue "D30" "1.479"
$ws.Range("E30").Value = "  -1.20%  "
Set-TextValue "D31" "0.05863"
$ws.Range("E31").Value = "  -5.97%  "
Set-TextValue "D32" "1.260"
$ws.Range("E32").Value = "  -2.73%  "
Set-TextValue "D33" "4.097"
$ws.Range("E33").Value = "  -1.45%  "
Set-TextValue "D34" "4.040"
$ws.Range("E34").Value = "  -1.76%  "
$ws.Range("E35").Value = "  -0.98%  "
Set-TextValue "D36" "1.169"
$ws.Range("E36").Value = "  -0.46%  "
Set-TextValue "D37" "0.7198"
$ws.Range("E37").Value = "  -1.44%  "
Set-TextValue "D38" "2.587"
$ws.Range("E38").Value = "  -0.87%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D39" "2.794"
$ws.Range("E39").Value = "  -1.09%  "
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").Value = "1.237.24"
$ws.Range("E40").Value = "  +2.04%  "
Set-TextValue "D41" "0.01776"
$ws.Range("E41").Value = "  -0.63%  "
Set-TextValue "D42" "0.9063"
$ws.Range("E42").Value = "  -1.89%  "
Set-TextValue "D43" "6.137"
$ws.Range("E43").Value = "  -2.56%  "
$ws.Range("D44").Value = "2.130.18"
$ws.Range("E44").Value = "  +5.30%  "
$ws.Range("E45").Value = "  -0.08%  "
Set-TextValue "D46" "101.93"
$ws.Range("E46").Value = "  -0.21%  "
Set-TextValue "D47" "67.03"
$ws.Range("E47").Value = "  +0.19%  "
Set-TextValue "D48" "7.300"
$ws.Range("E48").Value = "  +8.41%  "
$ws.Range("B49").Value = "BabyDogeCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue "D49" "0.00000000119"
$ws.Range("E49").Value = "  +0.62%  "
$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D50" "9.165"
$ws.Range("E50").Value = "  -0.73%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D51" "1.718"
$ws.Range("E51").Value = "  +2.82%  "
